$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") contains the date 45781 (2025-05-04) in rows 2-43.
# Update it to 45783 (2025-05-06) for every data row.
$ws.Range("C2:C43").Value = 45783
